# Commit: "Removed incorrect 'beta ratio' statement"
#
# The project-description bullet about the stock-calculator algorithm
# incorrectly referred to a "Beta Ratio" (beta is not a ratio, P/E is).
# Fix the wording from:
#   "...the company's P/E and Beta Ratio, Income Statement..."
# to:
#   "...the company's P/E ratio, Beta, Income Statement..."

$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "P/E and Beta Ratio",  # FindText
    $true,                  # MatchCase
    $true,                  # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $true,                  # Format
    "P/E ratio, Beta",      # ReplaceWith
    2                       # Replace (wdReplaceAll)
)
